$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the detailed-description text first (matches the order new
# shared strings were introduced), then the names/progress columns.
$ws.Range("E7").Value = "Class Creation done (see src) and manager in progress"
$ws.Range("E8").Value = "Class Creation done (see src) and manager in progress (group with Customer)"

# Row 7 - "Creating classes for Customers and its Manager"
$ws.Range("B7").Value = "Thomas, Yamid"
$ws.Range("D7").Value = "In Progress"

# Row 8 - "Creating classes for  Bills and its Manager(if required)"
$ws.Range("B8").Value = "Thomas, Yamid"
$ws.Range("D8").Value = "In Progress"

# Widen column E to fit the new, longer description text (69 characters
# wide in the saved file; the runtime adds ~5/6 of a character of padding
# when converting the ColumnWidth property to the stored width, so we back
# that padding out here)
$ws.Columns.Item(5).ColumnWidth = 68.16666666666667

# Move the active selection to A16, matching where the user clicked next
$ws.Range("A16").Select()
